# se modif data para regresion en pre prod R31
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the "regresion" data row (row 11) with the new pre-prod environment info
# (order matches the shared-strings table insertion order: Polizas, domain, URL, name)
$ws.Range("F11").Value = "Polizas"
$ws.Range("A11").Value = "i-preproducciongestion.segurossura.com.ar"
$ws.Range("B11").Value = "https://i-preproducciongestion.segurossura.com.ar/pc/PolicyCenter.do"
$ws.Range("J11").Value = "Baioni Alejandro Luis"
$ws.Range("E11").Value = "Prueba"
$ws.Range("G11").Value = 24741860
$ws.Range("I11").Value = 234
$ws.Range("N11").Value = 305

# Re-assigning .Value on A11 drops its "quote prefix" cell format (the
# Ambiente column is entered as quoted text), so restore the original
# formatting from a sibling cell that carries the same style.
$ws.Range("A6").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update selection to match the row that was edited (Ranorex selects the
# whole row after entering the last field of the record)
$ws.Activate()
$ws.Range("F11").Activate()
$ws.Rows("11:11").Select()
